# ---------------------------------------------------------------------------
# Edit script for assignment_05_exercise_09_JenaBinay.docx
# Applies the changes described in the commit "2020-10-03: Week 5 - Exercise
# 09 Student Survey":
#   1. Rewrite the covariance explanation paragraph.
#   2. Tighten the three variable-measurement bullet points.
#   3. Replace the "Correlation Coefficient" bullet list + "Coefficient of
#      Determination" discussion with a new prose paragraph (keeps both
#      R^2 equations, adds a bold run).
#   4. Rewrite the "Yes, Watching TV..." conclusion paragraph.
#   5. Replace the trailing partial-correlation bullet analysis with a
#      single new paragraph.
#   6. (numbering.xml num 1004/1005 definitions become orphaned once their
#      paragraphs are deleted - handled automatically, nothing to edit there
#      directly.)
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Find-ParagraphIndex($doc, $startsWith) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text
        if ($t.StartsWith($startsWith)) {
            return $i
        }
    }
    return -1
}

function Xml-Escape($text) {
    $escaped = $text -replace '&', '&amp;'
    $escaped = $escaped -replace '<', '&lt;'
    $escaped = $escaped -replace '>', '&gt;'
    return $escaped
}

function Build-RunXml($text, $bold) {
    $escaped = Xml-Escape $text
    if ($bold) {
        return '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
    } else {
        return '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
    }
}

# Exact copy of the R^2 "coefficient of determination" equation used twice
# in the document.
$omathR2 = '<m:oMath><m:sSup><m:e><m:r><m:t>R</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup></m:oMath>'

function Wrap-PkgXml($innerParagraphXml) {
    return '<?xml version="1.0"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math">' +
        '<w:body><w:p>' + $innerParagraphXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# Replace the run-content of a paragraph (identified by index) with the
# supplied inner-paragraph OOXML (runs / oMath), leaving the paragraph mark
# (and therefore its pPr / style) untouched.
function Set-ParagraphInnerXml($doc, $paraIndex, $innerXml) {
    $p = $doc.Paragraphs($paraIndex)
    $r = $p.Range
    $target = $doc.Range($r.Start, $r.End - 1)
    $full = Wrap-PkgXml $innerXml
    $target.InsertXML($full)
}

# ---------------------------------------------------------------------------
# 1. Covariance explanation paragraph (section a.)
# ---------------------------------------------------------------------------

$idx = Find-ParagraphIndex $d "Covariance is the simplest"
if ($idx -eq -1) { throw "Could not locate covariance paragraph" }

$inner =
    (Build-RunXml "Covariance is the simplest way to look or compare the two variables. It helps in understanding whether the two variables in question simultaneously i.e. co-vary with each other." $false) +
    (Build-RunXml " " $false) +
    (Build-RunXml "This metric is widely accepted as an indicator of related-ness amongst two variables." $false) +
    (Build-RunXml " " $false) +
    (Build-RunXml "A positive covariance indicates directional bias or movement in one of variables from the mean would also mean same directional deviation with the other variable w.r.t the eman. If this deviation is in opposite direction from the mean then the covariance is negative." $false)

Set-ParagraphInnerXml $d $idx $inner

# ---------------------------------------------------------------------------
# 2. Variable measurement bullet points (section b.)
# ---------------------------------------------------------------------------

$idx = Find-ParagraphIndex $d "TimeReading: This seems to be in hours"
if ($idx -eq -1) { throw "Could not locate TimeReading bullet" }
$inner = Build-RunXml "TimeReading: in hours" $false
Set-ParagraphInnerXml $d $idx $inner

$idx = Find-ParagraphIndex $d "TimeTV: This seems to be in minutes"
if ($idx -eq -1) { throw "Could not locate TimeTV bullet" }
$inner = Build-RunXml "TimeTV: in minutes" $false
Set-ParagraphInnerXml $d $idx $inner

$idx = Find-ParagraphIndex $d "Happiness: This seems either percentage"
if ($idx -eq -1) { throw "Could not locate Happiness bullet" }
$inner = Build-RunXml "Happiness: some numeric score could be percentages too" $false
Set-ParagraphInnerXml $d $idx $inner

# ---------------------------------------------------------------------------
# 3. Correlation Coefficient / Coefficient of Determination section
#    (section e.): old structure was
#       FirstParagraph: "Based on Correlation Coefficient..."
#       6x Compact bullets (numId 1004)
#       FirstParagraph: "However with Coefficient of Determination or the" + R^2
#       BodyText: "We can say that"
#       6x Compact bullets (numId 1005)
#    New structure is a single FirstParagraph with the two equations kept
#    (one mid-sentence, one new one appended near the end) plus a bold run.
# ---------------------------------------------------------------------------

$startIdx = Find-ParagraphIndex $d "Based on Correlation Coefficient, I already explained"
if ($startIdx -eq -1) { throw "Could not locate 'Based on Correlation Coefficient' paragraph" }

$endIdx = Find-ParagraphIndex $d "Happiness and Gender have covariability of 2.46%"
if ($endIdx -eq -1) { throw "Could not locate 'Happiness and Gender have covariability' paragraph" }

# Delete everything from the paragraph right after the opening paragraph
# through to (and including) the final covariability bullet, leaving just
# the opening "Based on Correlation Coefficient..." paragraph to rewrite.
$pFrom = $d.Paragraphs($startIdx + 1)
$pTo = $d.Paragraphs($endIdx)
$delRange = $d.Range($pFrom.Range.Start, $pTo.Range.End)
$delRange.Delete()

$inner =
    (Build-RunXml "Although we cannot make direct conclusions about causality from a correlation, we" $false) +
    (Build-RunXml " " $false) +
    (Build-RunXml "can take the correlation coefficient a step further by squaring it. The correlation" $false) +
    (Build-RunXml " " $false) +
    (Build-RunXml "coefficient squared (known as the coefficient of determination," $false) +
    (Build-RunXml " " $false) +
    $omathR2 +
    (Build-RunXml ") is a measure" $false) +
    (Build-RunXml " " $false) +
    (Build-RunXml "of the amount of variability in one variable that is shared by the other." $false) +
    (Build-RunXml " " $false) +
    (Build-RunXml "In our student survey example the correlation coefficient tells us that the watching TV is negatively related to reading. However we don’t know how much percent of affected reading time is because of watching TV. This is where" $false) +
    (Build-RunXml " " $false) +
    (Build-RunXml "Coefficient of Determination" $true) +
    (Build-RunXml " " $false) +
    (Build-RunXml "comes handy. It shows us what percent of reading is affected by watching TV. So above" $false) +
    (Build-RunXml " " $false) +
    $omathR2 +
    (Build-RunXml " " $false) +
    (Build-RunXml "matrix shows that the 77% of the time the reading is affected by watching TV." $false)

Set-ParagraphInnerXml $d $startIdx $inner

# ---------------------------------------------------------------------------
# 4. "Yes, Watching TV causes Students..." conclusion (section f.)
# ---------------------------------------------------------------------------

$idx = Find-ParagraphIndex $d "Yes, Watching TV causes Students to read less."
if ($idx -eq -1) { throw "Could not locate Watching TV conclusion paragraph" }

$inner = Build-RunXml "Yes, Watching TV causes Students to read less. Based on correlation test of student survey attributes we can say reading is affected by watching TV. Also as we have seen coefficient of determination also shows as much as 77% of the time reading time is affected by watching TV." $false

Set-ParagraphInnerXml $d $idx $inner

# ---------------------------------------------------------------------------
# 5. Partial correlation closing paragraphs (section g.): three paragraphs
#    collapse into one.
# ---------------------------------------------------------------------------

$startIdx = Find-ParagraphIndex $d "So here we see partial correlation coefficient is -0.886"
if ($startIdx -eq -1) { throw "Could not locate partial correlation paragraph" }

$endIdx = Find-ParagraphIndex $d "Hence we can conclude that the gender variable has least effect"
if ($endIdx -eq -1) { throw "Could not locate 'Hence we can conclude' paragraph" }

$pFrom = $d.Paragraphs($startIdx + 1)
$pTo = $d.Paragraphs($endIdx)
$delRange = $d.Range($pFrom.Range.Start, $pTo.Range.End)
$delRange.Delete()

$inner = Build-RunXml "Partial correlation analysis using TimeTv, TimeReading and Happiness shows that the time watching TV is negatively affecting reading time. Also when we keep Happiness constant doesn’t affect much the relation between watching TV and reading time. With correlation test we had r = -0.88 where as with partial test we get partial correlation of -0.87." $false

Set-ParagraphInnerXml $d $startIdx $inner

Write-Host "Done"
